$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    ,@(7, 23.2, 19.7, 28, 21.4, 21.6)
    ,@(8, 23.8, 24.9, 29.5, 24.2, 22.6)
    ,@(18, 25.8, 29.4, 28.7, 32.200000000000003, 29.4)
    ,@(21, 33.1, 26.1, 24.4, 22.5, 25.5)
    ,@(22, 32.9, 28.3, 25.7, 24.4, 25.5)
    ,@(23, 24.2, 26.7, 26.2, 28.8, 22.4)
    ,@(28, 28.7, 23.3, 23.2, 25.2, 41)
    ,@(31, 37.1, 24.3, 25.7, 23.3, 22)
    ,@(36, 31, 28.9, 27.3, 27.9, 25.6)
    ,@(46, 35.200000000000003, 22.1, 20.3, 20.9, 20.3)
    ,@(66, 59.2, 33.9, 35.9, 50.6, 40.6)
    ,@(67, 24.6, 21.4, 34.5, 32.299999999999997, 25)
    ,@(83, 30.2, 40.799999999999997, 37.4, 35.9, 33.200000000000003)
    ,@(93, 30.8, 33.5, 34.9, 34.4, 33.5)
    ,@(95, 41.4, 29.5, 36, 33, 29.6)
    ,@(96, 37.9, 39.9, 42.6, 41.8, 37.799999999999997)
    ,@(97, 36.799999999999997, 37.299999999999997, 37.9, 36.6, 32.6)
    ,@(300, 21.4, 22.4, 34.6, 21.9, 26.6)
    ,@(397, 61.4, 29.9, 27.8, 16.399999999999999, 28.5)
    ,@(403, 56.6, 21.6, 9.5, 21.2, 18.7)
    ,@(407, 53.6, 24.3, 21.8, 25, 20.8)
    ,@(413, 49.9, 32.6, 33.200000000000003, 31.8, 33.700000000000003)
    ,@(421, 50.6, 19.7, 20.100000000000001, 21.6, 22.9)
    ,@(428, 46.8, 21.2, 16.899999999999999, 22.2, 19.100000000000001)
    ,@(432, 48, 25.7, 27.4, 24, 22.8)
    ,@(439, 30.6, 21.9, 22, 26.2, 21.8)
    ,@(444, 16.8, 22.5, 17.8, 20.9, 19.3)
    ,@(445, 21.3, 23.4, 21.2, 27.6, 20.9)
    ,@(446, 28.2, 21.5, 22.8, 26.5, 26.1)
    ,@(457, 27.3, 29.8, 32.6, 26.2, 18.7)
    ,@(460, 35.1, 26, 23.2, 26.5, 23)
    ,@(465, 31, 24.7, 22.5, 25.8, 17.899999999999999)
    ,@(474, 27.3, 19.7, 20.8, 27.4, 21.5)
    ,@(475, 21.4, 20.7, 26, 19.899999999999999, 19.2)
    ,@(481, 30.5, 24.4, 25.7, 26, 24.4)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 16).Value = $entry[1]   # P
    $ws.Cells.Item($row, 17).Value = $entry[2]   # Q
    $ws.Cells.Item($row, 18).Value = $entry[3]   # R
    $ws.Cells.Item($row, 19).Value = $entry[4]   # S
    $ws.Cells.Item($row, 20).Value = $entry[5]   # T
}

$ws.Range("U481").Select()
